$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Area" column (B) from ITA17 to SLO, and "Station" column (C) from 45BIS to 45bis
# for all data rows (2 through 26).
$ws.Range("B2:B26").Value = "SLO"
$ws.Range("C2:C26").Value = "45bis"
